$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 51; this shifts existing rows 51-75 down to 52-76
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new weekly record.
# All descriptive columns mirror the surrounding rows for this market/product;
# only the date (D) and volume (J) differ for this new entry.
$ws.Range("A51").Value = 2
$ws.Range("B51").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 44825
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = 100112022
$ws.Range("G51").Value = "Arveja Verde"
$ws.Range("H51").Value = "Perfection"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 480
$ws.Range("K51").Value = 28000
$ws.Range("L51").Value = 30000
$ws.Range("M51").Value = 29000
$ws.Range("N51").Value = "$/malla 25 kilos"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 1160
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
